# Refactors corrected offers generation:
# - Rotates the promotional data (A, C, D, E columns) for rows 2-8 up by one
#   row (row 2's data wraps around to row 8), so the previously retained
#   promotional price from the consolidated file lines up with the correct
#   product row.
# - Renames the "Seção" (column B) values for clarity: the "ALTO GIRO" rows
#   are now split into "#01 ALTO GIRO" (rows 2-7) and "#02 ALTO GIRO"
#   (rows 8-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture current (pre-edit) values for the A/C/D/E columns of
#     the rotating block (rows 2-8) before overwriting anything.
$blockStart = 2
$blockEnd = 8
$rowCount = $blockEnd - $blockStart + 1

$colA = @{}
$colC = @{}
$colD = @{}
$colE = @{}

for ($r = $blockStart; $r -le $blockEnd; $r++) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value2
    $colC[$r] = $ws.Cells.Item($r, 3).Value2
    $colD[$r] = $ws.Cells.Item($r, 4).Value2
    $colE[$r] = $ws.Cells.Item($r, 5).Value2
}

# --- Step 2: write back the rotated values so that row N gets what used to
#     be in row N+1 (row 8 gets what used to be in row 2).
for ($i = 0; $i -lt $rowCount; $i++) {
    $destRow = $blockStart + $i
    $srcRow = $blockStart + (($i + 1) % $rowCount)

    $ws.Cells.Item($destRow, 1).Value = $colA[$srcRow]
    $ws.Cells.Item($destRow, 3).Value = $colC[$srcRow]
    $ws.Cells.Item($destRow, 4).Value = $colD[$srcRow]
    $ws.Cells.Item($destRow, 5).Value = $colE[$srcRow]
}

# --- Step 3: rename the "Seção" column (B) values for clarity.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 2).Value = "#01 MERCEARIA - #01 ALTO GIRO"
}
for ($r = 8; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = "#01 MERCEARIA - #02 ALTO GIRO"
}
